# 2021-06 QLD Outbreak Paths.xlsx — "Add files via upload"
#
# Semantic edit (derived from the OOXML diff):
#   Sheet1 ("Outbreak Paths" / Table1):
#     - C23 (Target) relabeled from "C3 w50" to "Q18 w50"
#     - Row 24 updated with a new case:
#         A24 Date           1-Jul-21 (44380) -> 2-Jul-21 (44381)
#         B24 Source         "1 m50" -> "Q17 m29"
#         C24 Target         "1 m50" -> "Q19"
#         D24 Cluster        "Unknown" -> "Queensland"
#         F24 Link Label     "Airport Worker, Lives in Carindale" -> "Carindale Greek Community Centre"
#         G24 Variant        "Unknown" -> "Alpha (B.1.1.7)"
#         H24 Isolated       "Not isolated" (unchanged)
#     - Active selection moves to B24
#
#   "Date Colours" sheet (Table Date_Colours):
#     - Colour Code gradient (B2:B15) recomputed to spread across one more row
#     - New row 16 gains a Colour Code value (the old final swatch, #9e57a6)
#     - Active selection on this sheet set to A2

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Date Colours")

# ---- "Date Colours" sheet: refreshed gradient -----------------------------
$ws2.Range("B2").Value = "#f9f3f9"
$ws2.Range("B3").Value = "#f3e8f3"
$ws2.Range("B4").Value = "#ecdced"
$ws2.Range("B5").Value = "#e6d1e7"
$ws2.Range("B6").Value = "#e0c6e1"
$ws2.Range("B7").Value = "#d9bbdb"
$ws2.Range("B8").Value = "#d3afd5"
$ws2.Range("B9").Value = "#cda4cf"
$ws2.Range("B10").Value = "#c699c9"
$ws2.Range("B11").Value = "#c08ec3"
$ws2.Range("B12").Value = "#b983bd"
$ws2.Range("B13").Value = "#b278b8"
$ws2.Range("B14").Value = "#ac6db2"
$ws2.Range("B15").Value = "#a562ac"
$ws2.Range("B16").Value = "#9e57a6"

# ---- Sheet1: "Outbreak Paths" --------------------------------------------

# Row 23: rename the target label
$ws1.Range("C23").Value = "Q18 w50"

# Row 24: replace the placeholder case with the new, real case data
$ws1.Range("A24").Value = 44381
$ws1.Range("B24").Value = "Q17 m29"
$ws1.Range("C24").Value = "Q19"
$ws1.Range("D24").Value = "Queensland"
$ws1.Range("F24").Value = "Carindale Greek Community Centre"
$ws1.Range("G24").Value = "Alpha (B.1.1.7)"
$ws1.Range("H24").Value = "Not isolated"

# Restore the recorded selections: "Date Colours" -> A2, then Sheet1 (the
# active tab) -> B24.
$ws2.Range("A2").Select()
$ws1.Activate()
$ws1.Range("B24").Select()
